$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (FAPs -> Ccl21b -> Ccr7 -> Resolving-Mac)
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("H2").Value = 0.6219589999999999
$ws.Range("I2").Value = 0.6398583988494134
$ws.Range("J2").Value = 0.6398583988494134
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.626356
$ws.Range("N2").Value = 7.879068
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.544495250468
$ws.Range("R2").Value = 4.900457254211999
$ws.Range("S2").Value = 0.6398583988494134
$ws.Range("T2").Value = 0.6398583988494134

# Row 3 becomes what used to be row 4's data (MuSCs -> Ccl21b -> Ccr7 -> Resolving-Mac)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.116689
$ws.Range("H3").Value = 0.350067
$ws.Range("I3").Value = 0.3601416011505865
$ws.Range("J3").Value = 0.3601416011505865
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.626356
$ws.Range("N3").Value = 7.879068
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.306466855284
$ws.Range("R3").Value = 2.758201697556
$ws.Range("S3").Value = 0.3601416011505865
$ws.Range("T3").Value = 0.3601416011505865

# Remove now-unused rows 4 and 5 entirely
$ws.Range("A4:T5").Delete()
